$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: semantic "type" row
$ws.Range("A2").Value = "iaest-measure:horas-trabajadas"
$ws.Range("B2").Value = "iaest-measure:personas-residentes-viviendas-familiares"
$ws.Range("C2").Value = "sdmx-dimension:refArea"
$ws.Range("D2").Value = "sdmx-dimension:refArea"
$ws.Range("E2").Value = "sdmx-dimension:refArea"
$ws.Range("F2").Value = "iaest-measure:sexo"
$ws.Range("G2").Value = "null"
$ws.Range("H2").Value = "null"

# Row 3: dim/medida row
$ws.Range("A3").Value = "medida"
$ws.Range("B3").Value = "medida"
$ws.Range("C3").Value = "dim"
$ws.Range("D3").Value = "dim"
$ws.Range("E3").Value = "dim"
$ws.Range("F3").Value = "medida"
$ws.Range("G3").Value = "null"
$ws.Range("H3").Value = "null"

# Row 4: data type row
$ws.Range("A4").Value = "xsd:int"
$ws.Range("B4").Value = "xsd:int"
$ws.Range("C4").Value = "URI-Municipio"
$ws.Range("D4").Value = "URI-Provincia"
$ws.Range("E4").Value = "URI-Comunidad"
$ws.Range("F4").Value = "xsd:int"
$ws.Range("G4").Value = "null"
$ws.Range("H4").Value = "null"

# Row 5: removed entirely (was mapping file references)
$ws.Range("A5:H5").Delete()
